$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws2 = $wb.Worksheets.Item("ARM")
$ws3 = $wb.Worksheets.Item("BSM")
$ws4 = $wb.Worksheets.Item("CRP")
$ws5 = $wb.Worksheets.Item("CUL")
$ws6 = $wb.Worksheets.Item("GSM")
$ws7 = $wb.Worksheets.Item("LTW")
$ws8 = $wb.Worksheets.Item("WVR")

# ALC row 33
$ws1.Range("H33").Value = 1260.7778
$ws1.Range("I33").Value = 140.83333
$ws1.Range("J33").Value = 3500.6667
$ws1.Range("K33").Value = 140.83333
$ws1.Range("L33").Value = 3500.6667
$ws1.Range("M33").Value = 88.16667000000001
$ws1.Range("N33").Value = -3958.6667

# ALC row 62
$ws1.Range("H62").Value = 5135.25
$ws1.Range("I62").Value = 3275.5
$ws1.Range("J62").Value = 6995.0
$ws1.Range("K62").Value = 3275.5
$ws1.Range("L62").Value = 6995.0
$ws1.Range("M62").Value = -2651.5
$ws1.Range("N62").Value = -8243.0

# ALC row 65
$ws1.Range("H65").Value = 5135.25
$ws1.Range("I65").Value = 3275.5
$ws1.Range("J65").Value = 6995.0
$ws1.Range("K65").Value = 16377.5
$ws1.Range("L65").Value = 34975.0
$ws1.Range("M65").Value = -13257.5
$ws1.Range("N65").Value = -41215.0

# ALC row 82
$ws1.Range("H82").Value = 5833.3335
$ws1.Range("I82").Value = 5833.3335
$ws1.Range("K82").Value = 17500.0005
$ws1.Range("M82").Value = -17094.0005

# ALC row 85
$ws1.Range("H85").Value = 5833.3335
$ws1.Range("I85").Value = 5833.3335
$ws1.Range("K85").Value = 17500.0005
$ws1.Range("M85").Value = -16096.0005

# ALC row 138
$ws1.Range("H138").Value = 4221.154
$ws1.Range("J138").Value = 5496.875
$ws1.Range("L138").Value = 16490.625
$ws1.Range("N138").Value = -26770.625

# ARM row 74
$ws2.Range("H74").Value = 957.0
$ws2.Range("I74").Value = 1000.0
$ws2.Range("K74").Value = 1000.0
$ws2.Range("M74").Value = -126.0

# ARM row 77
$ws2.Range("H77").Value = 957.0
$ws2.Range("I77").Value = 1000.0
$ws2.Range("K77").Value = 5000.0
$ws2.Range("M77").Value = -632.0

# ARM row 88
$ws2.Range("H88").Value = 3332.0
$ws2.Range("I88").Value = 2000.0
$ws2.Range("J88").Value = 3998.0
$ws2.Range("K88").Value = 2000.0
$ws2.Range("L88").Value = 3998.0
$ws2.Range("M88").Value = -1594.0
$ws2.Range("N88").Value = -4810.0

# ARM row 91
$ws2.Range("H91").Value = 3332.0
$ws2.Range("I91").Value = 2000.0
$ws2.Range("J91").Value = 3998.0
$ws2.Range("K91").Value = 2000.0
$ws2.Range("L91").Value = 3998.0
$ws2.Range("M91").Value = -596.0
$ws2.Range("N91").Value = -6806.0

# ARM row 104
$ws2.Range("H104").Value = 33306.0
$ws2.Range("J104").Value = 33306.0
$ws2.Range("L104").Value = 33306.0
$ws2.Range("N104").Value = -40294.0

# ARM row 109
$ws2.Range("H109").Value = 20377.0
$ws2.Range("J109").Value = 20377.0
$ws2.Range("L109").Value = 20377.0
$ws2.Range("N109").Value = -23151.0

# BSM row 20
$ws3.Range("H20").Value = 1611.8572
$ws3.Range("I20").Value = 1542.1111
$ws3.Range("J20").Value = 1737.4
$ws3.Range("K20").Value = 1542.1111
$ws3.Range("L20").Value = 1737.4
$ws3.Range("M20").Value = -1295.1111
$ws3.Range("N20").Value = -2231.4

# BSM row 22
$ws3.Range("H22").Value = 129.0
$ws3.Range("I22").Value = 111.25
$ws3.Range("J22").Value = 200.0
$ws3.Range("K22").Value = 111.25
$ws3.Range("L22").Value = 200.0
$ws3.Range("M22").Value = 61.75
$ws3.Range("N22").Value = -546.0

# BSM row 94
$ws3.Range("H94").Value = 1184.2632
$ws3.Range("I94").Value = 1184.2632
$ws3.Range("K94").Value = 1184.2632
$ws3.Range("M94").Value = -733.2632000000001

# CRP row 31
$ws4.Range("H31").Value = 1786.1428
$ws4.Range("I31").Value = 730.0909
$ws4.Range("J31").Value = 2947.8
$ws4.Range("K31").Value = 730.0909
$ws4.Range("L31").Value = 2947.8
$ws4.Range("M31").Value = -435.0909
$ws4.Range("N31").Value = -3537.8

# CRP row 34
$ws4.Range("H34").Value = 1786.1428
$ws4.Range("I34").Value = 730.0909
$ws4.Range("J34").Value = 2947.8
$ws4.Range("K34").Value = 730.0909
$ws4.Range("L34").Value = 2947.8
$ws4.Range("M34").Value = -528.0909
$ws4.Range("N34").Value = -3351.8

# CUL row 33
$ws5.Range("H33").Value = 21.0
$ws5.Range("I33").Value = 21.0
$ws5.Range("J33").Value = 0.0
$ws5.Range("K33").Value = 126.0
$ws5.Range("L33").Value = 0.0
$ws5.Range("M33").Value = 157.0
$ws5.Range("N33").ClearContents()

# GSM row 11
$ws6.Range("H11").Value = 7134727.0
$ws6.Range("I11").Value = 11444444.0
$ws6.Range("J11").Value = 1593661.4
$ws6.Range("K11").Value = 11444444.0
$ws6.Range("L11").Value = 1593661.4
$ws6.Range("M11").Value = -11444305.0
$ws6.Range("N11").Value = -1593939.4

# GSM row 14
$ws6.Range("H14").Value = 7333766.5
$ws6.Range("J14").Value = 2500649.8
$ws6.Range("L14").Value = 2500649.8
$ws6.Range("N14").Value = -2500985.8

# GSM row 80
$ws6.Range("H80").Value = 3374.5
$ws6.Range("I80").Value = 3249.4
$ws6.Range("J80").Value = 4000.0
$ws6.Range("K80").Value = 3249.4
$ws6.Range("L80").Value = 4000.0
$ws6.Range("M80").Value = -2251.4
$ws6.Range("N80").Value = -5996.0

# GSM row 83
$ws6.Range("H83").Value = 3374.5
$ws6.Range("I83").Value = 3249.4
$ws6.Range("J83").Value = 4000.0
$ws6.Range("K83").Value = 16247.0
$ws6.Range("L83").Value = 20000.0
$ws6.Range("M83").Value = -11255.0
$ws6.Range("N83").Value = -29984.0

# GSM row 97
$ws6.Range("H97").Value = 0.0
$ws6.Range("I97").Value = 0.0
$ws6.Range("K97").Value = 0.0
$ws6.Range("M97").ClearContents()

# GSM row 102
$ws6.Range("H102").Value = 2468.75
$ws6.Range("I102").Value = 2623.2727
$ws6.Range("K102").Value = 2623.2727
$ws6.Range("M102").Value = -1001.2727

# GSM row 113
$ws6.Range("H113").Value = 0.0
$ws6.Range("I113").Value = 0.0
$ws6.Range("J113").Value = 0.0
$ws6.Range("K113").Value = 0.0
$ws6.Range("L113").Value = 0.0
$ws6.Range("M113").ClearContents()
$ws6.Range("N113").ClearContents()

# LTW row 22
$ws7.Range("H22").Value = 3257.8572
$ws7.Range("I22").Value = 0.0
$ws7.Range("J22").Value = 3257.8572
$ws7.Range("K22").Value = 0.0
$ws7.Range("L22").Value = 3257.8572
$ws7.Range("M22").ClearContents()
$ws7.Range("N22").Value = -3847.8572

# LTW row 27
$ws7.Range("H27").Value = 3257.8572
$ws7.Range("I27").Value = 0.0
$ws7.Range("J27").Value = 3257.8572
$ws7.Range("K27").Value = 0.0
$ws7.Range("L27").Value = 3257.8572
$ws7.Range("M27").ClearContents()
$ws7.Range("N27").Value = -3471.8572

# LTW row 55
$ws7.Range("H55").Value = 5627.25
$ws7.Range("I55").Value = 0.0
$ws7.Range("J55").Value = 5627.25
$ws7.Range("K55").Value = 0.0
$ws7.Range("L55").Value = 5627.25
$ws7.Range("M55").ClearContents()
$ws7.Range("N55").Value = -5973.25

# LTW row 68
$ws7.Range("H68").Value = 4218.1816
$ws7.Range("I68").Value = 2666.6667
$ws7.Range("J68").Value = 11200.0
$ws7.Range("K68").Value = 2666.6667
$ws7.Range("L68").Value = 11200.0
$ws7.Range("M68").Value = -1917.6667
$ws7.Range("N68").Value = -12698.0

# LTW row 71
$ws7.Range("H71").Value = 4218.1816
$ws7.Range("I71").Value = 2666.6667
$ws7.Range("J71").Value = 11200.0
$ws7.Range("K71").Value = 13333.3335
$ws7.Range("L71").Value = 56000.0
$ws7.Range("M71").Value = -9589.3335
$ws7.Range("N71").Value = -63488.0

# WVR row 17
$ws8.Range("H17").Value = 4.0
$ws8.Range("I17").Value = 4.0
$ws8.Range("K17").Value = 4.0
$ws8.Range("M17").Value = 168.0

# WVR row 54
$ws8.Range("H54").Value = 8441.4
$ws8.Range("I54").Value = 8441.4
$ws8.Range("K54").Value = 8441.4
$ws8.Range("M54").Value = -7921.4

# WVR row 100
$ws8.Range("H100").Value = 649.5
$ws8.Range("I100").Value = 500.0
$ws8.Range("J100").Value = 799.0
$ws8.Range("K100").Value = 1000.0
$ws8.Range("L100").Value = 1598.0
$ws8.Range("M100").Value = -459.0
$ws8.Range("N100").Value = -2680.0

# WVR row 122
$ws8.Range("H122").Value = 3099.0715
$ws8.Range("I122").Value = 3708.8
$ws8.Range("J122").Value = 1574.75
$ws8.Range("K122").Value = 11126.4
$ws8.Range("L122").Value = 4724.25
$ws8.Range("M122").Value = -8676.400000000001
$ws8.Range("N122").Value = -9624.25

# WVR row 132
$ws8.Range("H132").Value = 2654.2222
$ws8.Range("I132").Value = 2714.5715
$ws8.Range("J132").Value = 2443.0
$ws8.Range("K132").Value = 8143.7145
$ws8.Range("L132").Value = 7329.0
$ws8.Range("M132").Value = -5613.7145
$ws8.Range("N132").Value = -12389.0
